$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.698.35"
$ws.Range("E2").Value = "  -5.45%  "
$ws.Range("D3").Value = "2.573.10"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'301.84"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'96.47"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").Value = "'0.581"
$ws.Range("E7").Value = "  -3.20%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.566"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "'37.12"
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("D11").Value = "'0.0818"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("D12").Value = "'7.84"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "2.966.85"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "2.578.05"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "'0.893"
$ws.Range("E16").Value = "  -2.57%  "
$ws.Range("D17").Value = "'14.38"
$ws.Range("E17").Value = "  -3.67%  "
$ws.Range("D18").Value = "43.767.86"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("D19").Value = "'6.72"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "0.0₃0986"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'12.53"
$ws.Range("E21").Value = "  -2.86%  "
$ws.Range("D22").Value = "'73.44"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").Value = "'266.13"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "'29.38"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'10.28"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'38.36"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "'6.20"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").Value = "'3.59"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'2.21"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "'152.55"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "'2.79"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").Value = "'0.0816"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("E38").Value = "  -1.56%  "
$ws.Range("D39").Value = "'24.03"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "'16.78"
$ws.Range("E40").Value = "  +6.23%  "
$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("D43").Value = "'3.88"
$ws.Range("E43").Value = "  -4.21%  "
$ws.Range("D44").Value = "2.037.73"
$ws.Range("E44").Value = "  -4.72%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'87.95"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("E48").Value = "  +5.30%  "
$ws.Range("D49").Value = "2.831.08"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "'105.60"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("E51").Value = "  -3.52%  "
